{"js": "const replacements = [\n  [\"2024-02-05 Monday\", \"2024-02-06 Tuesday\"],\n  [\"606\u00d79=\", \"269\u00d77=\"],\n  [\"625\u00d79=\", \"558\u00d74=\"],\n  [\"789\u00d75=\", \"919\u00d78=\"],\n  [\"424\u00d73=\", \"941\u00d78=\"],\n  [\"915\u00d72=\", \"728\u00d74=\"],\n  [\"311\u00d72=\", \"949\u00d78=\"],\n  [\"544\u00d76=\", \"419\u00d79=\"],\n  [\"313\u00d75=\", \"897\u00d78=\"],\n  [\"989\u00d78=\", \"587\u00d78=\"],\n  [\"866\u00d77=\", \"114\u00d78=\"],\n  [\"665\u00d72=\", \"563\u00d76=\"],\n  [\"193\u00d78=\", \"672\u00d74=\"],\n  [\"695\u00d79=\", \"391\u00d75=\"],\n  [\"278\u00d77=\", \"328\u00d75=\"],\n  [\"895\u00d78=\", \"523\u00d78=\"],\n  [\"506\u00d77=\", \"977\u00d78=\"],\n  [\"764\u00d73=\", \"352\u00d76=\"],\n  [\"477\u00d72=\", \"232\u00d75=\"],\n  [\"605\u00d74=\", \"359\u00d75=\"],\n  [\"766\u00d77=\", \"600\u00d74=\"],\n  [\"949\u00d73=\", \"219\u00d79=\"],\n  [\"399\u00d79=\", \"735\u00d79=\"],\n  [\"693\u00d74=\", \"918\u00d73=\"],\n  [\"462\u00d78=\", \"114\u00d78=\"],\n  [\"816\u00d79=\", \"296\u00d75=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + from);\n  }\n  for (const r of results.items) {\n    r.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-05 Monday\", \"2024-02-06 Tuesday\"),\n    @(\"606\u00d79=\", \"269\u00d77=\"),\n    @(\"625\u00d79=\", \"558\u00d74=\"),\n    @(\"789\u00d75=\", \"919\u00d78=\"),\n    @(\"424\u00d73=\", \"941\u00d78=\"),\n    @(\"915\u00d72=\", \"728\u00d74=\"),\n    @(\"311\u00d72=\", \"949\u00d78=\"),\n    @(\"544\u00d76=\", \"419\u00d79=\"),\n    @(\"313\u00d75=\", \"897\u00d78=\"),\n    @(\"989\u00d78=\", \"587\u00d78=\"),\n    @(\"866\u00d77=\", \"114\u00d78=\"),\n    @(\"665\u00d72=\", \"563\u00d76=\"),\n    @(\"193\u00d78=\", \"672\u00d74=\"),\n    @(\"695\u00d79=\", \"391\u00d75=\"),\n    @(\"278\u00d77=\", \"328\u00d75=\"),\n    @(\"895\u00d78=\", \"523\u00d78=\"),\n    @(\"506\u00d77=\", \"977\u00d78=\"),\n    @(\"764\u00d73=\", \"352\u00d76=\"),\n    @(\"477\u00d72=\", \"232\u00d75=\"),\n    @(\"605\u00d74=\", \"359\u00d75=\"),\n    @(\"766\u00d77=\", \"600\u00d74=\"),\n    @(\"949\u00d73=\", \"219\u00d79=\"),\n    @(\"399\u00d79=\", \"735\u00d79=\"),\n    @(\"693\u00d74=\", \"918\u00d73=\"),\n    @(\"462\u00d78=\", \"114\u00d78=\"),\n    @(\"816\u00d79=\", \"296\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
